# ---------------------------------------------------------------------------
# OralPresentation.pptx edit script
#
# Summary of the change:
#  1. Slide "Before Earthquake - Early warning system" (slide 2) title gets
#     " in India" appended.
#  2. Three brand-new slides are inserted describing: early warning systems
#     in Japan, crowdsourced-phone research, and "After Earthquake" response.
#  3. The existing "During Earthquake" slide gains body content (measuring
#     earthquake size) plus a source textbox, and moves later in the deck.
#  4. Final slide order (by original/new p:sldId):
#       256, 257, 258, 261(Japan), 260(Research), 259(During EQ), 262(After EQ)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$enDash = [char]0x2013

# ---------------------------------------------------------------------------
# 1. "Before Earthquake - Early warning system" -> "... in India"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$titleRange2 = $s2.Shapes.Item(1).TextFrame.TextRange
$fullTitle2 = $titleRange2.Text
$needle = " Early warning system"
$idx2 = $fullTitle2.IndexOf($needle)
$run2 = $titleRange2.Characters($idx2 + 1, $needle.Length)
$run2.Text = " Early warning system in India"

# ---------------------------------------------------------------------------
# 2. Build the three new slides by duplicating the existing "During
#    Earthquake" slide (slide 4) so the new slides inherit its layout /
#    placeholder geometry, then overwrite their text.  Chain-duplicating
#    keeps the new slides in ascending id / position order right after the
#    original.
# ---------------------------------------------------------------------------
$sDuring = $p.Slides.Item(4)
$dupResearch = $sDuring.Duplicate()        # new id 260 -> "Research" slide
$dupJapan = $dupResearch.Duplicate()       # new id 261 -> "Japan" slide
$dupAfter = $dupJapan.Duplicate()          # new id 262 -> "After Earthquake" slide

# Current order: [.., During EQ(259), Research(260), Japan(261), After(262)]
# Desired order: [.., Japan(261), Research(260), During EQ(259), After(262)]
$sDuring.MoveTo(6)
$dupResearch.MoveTo(5)

# Re-fetch slides by position now that the order is final.
$sJapan = $p.Slides.Item(4)
$sResearch = $p.Slides.Item(5)
$sDuringEQ = $p.Slides.Item(6)
$sAfter = $p.Slides.Item(7)

# ---------------------------------------------------------------------------
# 2a. Slide 4 - "Before Earthquake - Early warning system in Japan"
# ---------------------------------------------------------------------------
$titleJapan = $sJapan.Shapes.Item(1).TextFrame.TextRange
$titleJapan.Text = "Before Earthquake " + $enDash + " Early warning system in Japan"

$bodyJapan = $sJapan.Shapes.Item(2)
$bodyJapan.TextFrame.TextRange.Text = "https://www.technologyreview.com/s/423279/how-japans-earthquake-and-tsunami-warning-systems-work/" + [char]13 + "http://www.jma.go.jp/jma/en/Activities/eew1.html" + [char]13
$bjPara = $bodyJapan.TextFrame.TextRange.Paragraphs()
$bjPara.Item(1).ActionSettings(1).Hyperlink.Address = "https://www.technologyreview.com/s/423279/how-japans-earthquake-and-tsunami-warning-systems-work/"
$bjPara.Item(2).ActionSettings(1).Hyperlink.Address = "http://www.jma.go.jp/jma/en/Activities/eew1.html"

# ---------------------------------------------------------------------------
# 2b. Slide 5 - "Before Earthquake - Early Warning system Research"
# ---------------------------------------------------------------------------
$titleResearch = $sResearch.Shapes.Item(1).TextFrame.TextRange
$titleResearch.Text = "Before Earthquake " + $enDash + " Early Warning system Research"

$bodyResearch = $sResearch.Shapes.Item(2)
$bodyResearch.Left = 1103312 / 12700
$bodyResearch.Top = 2052919 / 12700
$bodyResearch.Width = 8946541 / 12700
$bodyResearch.Height = 3148474 / 12700
$bodyResearch.TextFrame.TextRange.Text = "Crowdsourcing Mobile phones " + $enDash + " making use of GPS" + [char]13

$tbResearch = $sResearch.Shapes.AddTextbox(1, 1103312 / 12700, 6092042 / 12700, 8946541 / 12700, 246221 / 12700)
$tbResearch.TextFrame.WordWrap = -1
$tbResearch.TextFrame.AutoSize = 1
$tbResearch.Fill.Visible = 0
$tbResearch.TextFrame.TextRange.Text = "Source: http://www.sciencemag.org/news/2015/04/smart-phones-could-be-used-detect-earthquakes"
$tbResearch.TextFrame.TextRange.Font.Size = 10

# ---------------------------------------------------------------------------
# 2c. Slide 6 - "During Earthquake" (existing slide, now gains body content)
# ---------------------------------------------------------------------------
$bodyDuring = $sDuringEQ.Shapes.Item(2)
$bodyDuring.Left = 1103312 / 12700
$bodyDuring.Top = 2052918 / 12700
$bodyDuring.Width = 8946541 / 12700
$bodyDuring.Height = 2530957 / 12700
$bodyDuring.TextFrame.TextRange.Text = "Measuring size of earthquake" + [char]13 + "Magnitude" + [char]13 + "Energy" + [char]13 + "Intensity" + [char]13 + "Depth" + [char]13
$bdParas = $bodyDuring.TextFrame.TextRange.Paragraphs()
for ($i = 2; $i -le $bdParas.Count; $i++) {
  $bdParas.Item($i).IndentLevel = 2
}

$tbDuring = $sDuringEQ.Shapes.AddTextbox(1, 1187532 / 12700, 5902036 / 12700, 9001497 / 12700, 246221 / 12700)
$tbDuring.TextFrame.WordWrap = -1
$tbDuring.TextFrame.AutoSize = 1
$tbDuring.Fill.Visible = 0
$tbDuring.TextFrame.TextRange.Text = "Source: https://earthquake.usgs.gov/learn/topics/measure.php"
$tbDuring.TextFrame.TextRange.Font.Size = 10

# ---------------------------------------------------------------------------
# 2d. Slide 7 - "After Earthquake"
# ---------------------------------------------------------------------------
$titleAfter = $sAfter.Shapes.Item(1).TextFrame.TextRange
$titleAfter.Text = "After Earthquake"

$bodyAfter = $sAfter.Shapes.Item(2)
$bodyAfter.Left = 1103312 / 12700
$bodyAfter.Top = 2052918 / 12700
$bodyAfter.Width = 8946541 / 12700
$bodyAfter.Height = 2566583 / 12700
$bodyAfter.TextFrame.TextRange.Text = "Detection and alert for earthquake related activities" + [char]13 + "In Australia " + $enDash + " Geoscience Australia" + [char]13 + "monitors seismic data worldwide" + [char]13 + "Analyzed automatically and reviewed within 10 minutes of origin time" + [char]13 + "Results in alerts like tsunami alerts" + [char]13 + "Earthquake database" + [char]13
$baParas = $bodyAfter.TextFrame.TextRange.Paragraphs()
$baParas.Item(3).IndentLevel = 2
$baParas.Item(4).IndentLevel = 2
$baParas.Item(5).IndentLevel = 2

$tbAfter = $sAfter.Shapes.AddTextbox(1, 1199408 / 12700, 5510151 / 12700, 9084623 / 12700, 400110 / 12700)
$tbAfter.TextFrame.WordWrap = -1
$tbAfter.TextFrame.AutoSize = 1
$tbAfter.Fill.Visible = 0
$tbAfter.TextFrame.TextRange.Text = "Source: http://www.ga.gov.au/scientific-topics/hazards/earthquake/capabilties/monitoring" + [char]13 + "Source: http://www.ga.gov.au/earthquakes/searchQuake.do"
$tbAfter.TextFrame.TextRange.Font.Size = 10
$tbAfterParas = $tbAfter.TextFrame.TextRange.Paragraphs()
$tbAfterParas.Item(1).ActionSettings(1).Hyperlink.Address = "http://www.ga.gov.au/scientific-topics/hazards/earthquake/capabilties/monitoring"
